$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-07 22:18:20"
$ws.Range("O2").Value = "-0.7 °C"
$ws.Range("E3").Value = "2026-02-07 22:18:23"
$ws.Range("E4").Value = "2026-02-07 22:18:25"
$ws.Range("E5").Value = "2026-02-07 22:18:28"
$ws.Range("E6").Value = "2026-02-07 22:18:30"
$ws.Range("O6").Value = "10.6 °C"
$ws.Range("E7").Value = "2026-02-07 22:18:32"
$ws.Range("E8").Value = "2026-02-07 22:18:35"
$ws.Range("J8").Value = "1005.2 hPa"
$ws.Range("E9").Value = "2026-02-07 22:18:37"
$ws.Range("N9").Value = "4.0 °C 21:40 TU"
$ws.Range("O9").Value = "10.7 °C"
$ws.Range("E10").Value = "2026-02-07 22:18:40"
$ws.Range("H10").Value = "'77%"
$ws.Range("O10").Value = "8.2 °C"
$ws.Range("E11").Value = "2026-02-07 22:18:42"
$ws.Range("O11").Value = "3.4 °C"
$ws.Range("E12").Value = "2026-02-07 22:18:45"
$ws.Range("O12").Value = "10.1 °C"
$ws.Range("E13").Value = "2026-02-07 22:18:47"
$ws.Range("J13").Value = "1006.5 hPa"
$ws.Range("O13").Value = "3.0 °C"
$ws.Range("E14").Value = "2026-02-07 22:18:49"
$ws.Range("E15").Value = "2026-02-07 22:18:51"
$ws.Range("N15").Value = "4.4 °C 21:59 TU"
$ws.Range("O15").Value = "10.4 °C"
$ws.Range("E16").Value = "2026-02-07 22:18:53"
$ws.Range("H16").Value = "'60%"
$ws.Range("E17").Value = "2026-02-07 22:18:56"
$ws.Range("L17").Value = "74.5 km/h - 221º 21:38 TU"
$ws.Range("E18").Value = "2026-02-07 22:18:58"
$ws.Range("O18").Value = "9.8 °C"
$ws.Range("E19").Value = "2026-02-07 22:19:01"
$ws.Range("O19").Value = "3.9 °C"
$ws.Range("E20").Value = "2026-02-07 22:19:03"
$ws.Range("I20").Value = "1.5 mm"
$ws.Range("E21").Value = "2026-02-07 22:19:06"
$ws.Range("O21").Value = "4.9 °C"
$ws.Range("E22").Value = "2026-02-07 22:19:08"
$ws.Range("I22").Value = "2.0 mm"
$ws.Range("L22").Value = "34.6 km/h - 114º 21:54 TU"
$ws.Range("O22").Value = "-6.3 °C"
$ws.Range("E23").Value = "2026-02-07 22:19:11"
$ws.Range("E24").Value = "2026-02-07 22:19:13"
$ws.Range("I24").Value = "1.0 mm"
$ws.Range("E25").Value = "2026-02-07 22:19:16"
$ws.Range("H25").Value = "'78%"
$ws.Range("E26").Value = "2026-02-07 22:19:18"
$ws.Range("H26").Value = "'72%"
$ws.Range("L26").Value = "33.5 km/h - 224º 21:59 TU"
$ws.Range("E27").Value = "2026-02-07 22:19:21"
$ws.Range("H27").Value = "'83%"
$ws.Range("L27").Value = "34.9 km/h - 275º 21:10 TU"
$ws.Range("O27").Value = "-4.1 °C"
$ws.Range("E28").Value = "2026-02-07 22:19:23"
$ws.Range("L28").Value = "30.2 km/h - 221º 21:55 TU"
$ws.Range("E29").Value = "2026-02-07 22:19:26"
$ws.Range("E30").Value = "2026-02-07 22:19:28"
$ws.Range("O30").Value = "9.7 °C"
$ws.Range("E31").Value = "2026-02-07 22:19:31"
$ws.Range("N31").Value = "9.4 °C 21:43 TU"
$ws.Range("O31").Value = "11.1 °C"
$ws.Range("E32").Value = "2026-02-07 22:19:33"
$ws.Range("I32").Value = "0.7 mm"
$ws.Range("E33").Value = "2026-02-07 22:19:36"
$ws.Range("H33").Value = "'84%"
$ws.Range("E34").Value = "2026-02-07 22:19:38"
$ws.Range("O34").Value = "-2.3 °C"
$ws.Range("E35").Value = "2026-02-07 22:19:40"
$ws.Range("E36").Value = "2026-02-07 22:19:43"
$ws.Range("J36").Value = "1004.1 hPa"
$ws.Range("O36").Value = "11.3 °C"
$ws.Range("E37").Value = "2026-02-07 22:19:45"
$ws.Range("E38").Value = "2026-02-07 22:19:48"
$ws.Range("E39").Value = "2026-02-07 22:19:50"
$ws.Range("H39").Value = "'66%"
$ws.Range("E40").Value = "2026-02-07 22:19:53"
$ws.Range("E41").Value = "2026-02-07 22:19:55"
$ws.Range("E42").Value = "2026-02-07 22:19:58"
$ws.Range("E43").Value = "2026-02-07 22:20:00"
$ws.Range("O43").Value = "7.6 °C"
$ws.Range("E44").Value = "2026-02-07 22:20:03"
$ws.Range("E45").Value = "2026-02-07 22:20:05"
$ws.Range("H45").Value = "'60%"
$ws.Range("O45").Value = "4.2 °C"
$ws.Range("E46").Value = "2026-02-07 22:20:08"
$ws.Range("O46").Value = "9.2 °C"
